$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.293.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.43%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.566.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.91%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'618.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +2.53%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +3.95%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.564.25"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.88%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.00%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +2.31%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +5.35%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'7.41"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +6.81%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +3.96%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'33.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +5.50%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +1.04%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.169.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.91%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.568.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.86%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'68.359.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.66%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  -0.12%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +6.13%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'16.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +6.67%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'10.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +12.01%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'454.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.87%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +4.10%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'78.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.60%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +2.53%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'3.708.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.90%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -0.14%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'9.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +12.24%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +4.08%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +10.67%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +3.38%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E33").Value = "'  -0.03%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'6.39"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +4.49%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'26.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.97%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +4.63%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.559.69"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +2.05%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +3.50%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +8.94%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D41").Value = "'181.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +4.32%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.0918"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +5.01%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.04%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'5.60"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +3.69%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +11.78%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +2.07%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'46.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.94%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +3.89%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'2.66"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +4.38%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'7.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +3.59%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +7.71%  "
$ws.Range("E51").Style = "Normal"
